$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dag1"
$ws.Range("C2").Value = "Lama4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 18.76192366666666
$ws.Range("H2").Value = 56.285771
$ws.Range("I2").Value = 0.1222461152048115
$ws.Range("J2").Value = 0.1222461152048115
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 79.45695600000001
$ws.Range("N2").Value = 238.370868
$ws.Range("O2").Value = 0.3098915743290338
$ws.Range("P2").Value = 0.3098915743290338
$ws.Range("Q2").Value = 1490.765343257692
$ws.Range("R2").Value = 13416.88808931923
$ws.Range("S2").Value = 0.03788304109642748
$ws.Range("T2").Value = 0.03788304109642748

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dag1"
$ws.Range("C3").Value = "Lama4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 18.76192366666666
$ws.Range("H3").Value = 56.285771
$ws.Range("I3").Value = 0.1222461152048115
$ws.Range("J3").Value = 0.1222461152048115
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 164.182683
$ws.Range("N3").Value = 492.548049
$ws.Range("O3").Value = 0.6403319819152736
$ws.Range("P3").Value = 0.6403319819152735
$ws.Range("Q3").Value = 3080.382965834531
$ws.Range("R3").Value = 27723.44669251078
$ws.Range("S3").Value = 0.07827809723053983
$ws.Range("T3").Value = 0.07827809723053981

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Dag1"
$ws.Range("C4").Value = "Lama4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 18.76192366666666
$ws.Range("H4").Value = 56.285771
$ws.Range("I4").Value = 0.1222461152048115
$ws.Range("J4").Value = 0.1222461152048115
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.07725
$ws.Range("N4").Value = 0.23175
$ws.Range("O4").Value = 0.0003012841835637129
$ws.Range("P4").Value = 0.0003012841835637129
$ws.Range("Q4").Value = 1.44935860325
$ws.Range("R4").Value = 13.04422742925
$ws.Range("S4").Value = 0.00003683082101331723
$ws.Range("T4").Value = 0.00003683082101331723

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Dag1"
$ws.Range("C5").Value = "Lama4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 18.76192366666666
$ws.Range("H5").Value = 56.285771
$ws.Range("I5").Value = 0.1222461152048115
$ws.Range("J5").Value = 0.1222461152048115
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.68555166666667
$ws.Range("N5").Value = 38.056655
$ws.Range("O5").Value = 0.04947515957212899
$ws.Range("P5").Value = 0.04947515957212898
$ws.Range("Q5").Value = 238.0053520395561
$ws.Range("R5").Value = 2142.048168356005
$ws.Range("S5").Value = 0.006048146056830915
$ws.Range("T5").Value = 0.006048146056830913

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Dag1"
$ws.Range("C6").Value = "Lama4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 38.59812166666666
$ws.Range("H6").Value = 115.794365
$ws.Range("I6").Value = 0.2514918252404857
$ws.Range("J6").Value = 0.2514918252404857
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 79.45695600000001
$ws.Range("N6").Value = 238.370868
$ws.Range("O6").Value = 0.3098915743290338
$ws.Range("P6").Value = 0.3098915743290338
$ws.Range("Q6").Value = 3066.88925495098
$ws.Range("R6").Value = 27602.00329455882
$ws.Range("S6").Value = 0.07793519765465635
$ws.Range("T6").Value = 0.07793519765465635

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Dag1"
$ws.Range("C7").Value = "Lama4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 38.59812166666666
$ws.Range("H7").Value = 115.794365
$ws.Range("I7").Value = 0.2514918252404857
$ws.Range("J7").Value = 0.2514918252404857
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 164.182683
$ws.Range("N7").Value = 492.548049
$ws.Range("O7").Value = 0.6403319819152736
$ws.Range("P7").Value = 0.6403319819152735
$ws.Range("Q7").Value = 6337.143173993764
$ws.Range("R7").Value = 57034.28856594388
$ws.Range("S7").Value = 0.1610382588917298
$ws.Range("T7").Value = 0.1610382588917298

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Dag1"
$ws.Range("C8").Value = "Lama4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 38.59812166666666
$ws.Range("H8").Value = 115.794365
$ws.Range("I8").Value = 0.2514918252404857
$ws.Range("J8").Value = 0.2514918252404857
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.07725
$ws.Range("N8").Value = 0.23175
$ws.Range("O8").Value = 0.0003012841835637129
$ws.Range("P8").Value = 0.0003012841835637129
$ws.Range("Q8").Value = 2.98170489875
$ws.Range("R8").Value = 26.83534408875
$ws.Range("S8").Value = 0.0000757705092405277
$ws.Range("T8").Value = 0.0000757705092405277

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Dag1"
$ws.Range("C9").Value = "Lama4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 38.59812166666666
$ws.Range("H9").Value = 115.794365
$ws.Range("I9").Value = 0.2514918252404857
$ws.Range("J9").Value = 0.2514918252404857
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 12.68555166666667
$ws.Range("N9").Value = 38.056655
$ws.Range("O9").Value = 0.04947515957212899
$ws.Range("P9").Value = 0.04947515957212898
$ws.Range("Q9").Value = 489.6384666387861
$ws.Range("R9").Value = 4406.746199749075
$ws.Range("S9").Value = 0.01244259818485901
$ws.Range("T9").Value = 0.012442598184859

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Dag1"
$ws.Range("C10").Value = "Lama4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.691683333333334
$ws.Range("H10").Value = 20.07505
$ws.Range("I10").Value = 0.04360066196912097
$ws.Range("J10").Value = 0.04360066196912097
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 79.45695600000001
$ws.Range("N10").Value = 238.370868
$ws.Range("O10").Value = 0.3098915743290338
$ws.Range("P10").Value = 0.3098915743290338
$ws.Range("Q10").Value = 531.7007881826
$ws.Range("R10").Value = 4785.307093643401
$ws.Range("S10").Value = 0.01351147777939893
$ws.Range("T10").Value = 0.01351147777939893

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Dag1"
$ws.Range("C11").Value = "Lama4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 6.691683333333334
$ws.Range("H11").Value = 20.07505
$ws.Range("I11").Value = 0.04360066196912097
$ws.Range("J11").Value = 0.04360066196912097
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 164.182683
$ws.Range("N11").Value = 492.548049
$ws.Range("O11").Value = 0.6403319819152736
$ws.Range("P11").Value = 0.6403319819152735
$ws.Range("Q11").Value = 1098.65852345305
$ws.Range("R11").Value = 9887.926711077451
$ws.Range("S11").Value = 0.02791889829150513
$ws.Range("T11").Value = 0.02791889829150512

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Dag1"
$ws.Range("C12").Value = "Lama4"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 6.691683333333334
$ws.Range("H12").Value = 20.07505
$ws.Range("I12").Value = 0.04360066196912097
$ws.Range("J12").Value = 0.04360066196912097
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.07725
$ws.Range("N12").Value = 0.23175
$ws.Range("O12").Value = 0.0003012841835637129
$ws.Range("P12").Value = 0.0003012841835637129
$ws.Range("Q12").Value = 0.5169325375
$ws.Range("R12").Value = 4.652392837500001
$ws.Range("S12").Value = 0.00001313618984420404
$ws.Range("T12").Value = 0.00001313618984420404

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Dag1"
$ws.Range("C13").Value = "Lama4"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 6.691683333333334
$ws.Range("H13").Value = 20.07505
$ws.Range("I13").Value = 0.04360066196912097
$ws.Range("J13").Value = 0.04360066196912097
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 12.68555166666667
$ws.Range("N13").Value = 38.056655
$ws.Range("O13").Value = 0.04947515957212899
$ws.Range("P13").Value = 0.04947515957212898
$ws.Range("Q13").Value = 84.88769466197223
$ws.Range("R13").Value = 763.9892519577501
$ws.Range("S13").Value = 0.002157149708372716
$ws.Range("T13").Value = 0.002157149708372715

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Dag1"
$ws.Range("C14").Value = "Lama4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 89.42491666666666
$ws.Range("H14").Value = 268.27475
$ws.Range("I14").Value = 0.5826613975855819
$ws.Range("J14").Value = 0.5826613975855818
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 79.45695600000001
$ws.Range("N14").Value = 238.370868
$ws.Range("O14").Value = 0.3098915743290338
$ws.Range("P14").Value = 0.3098915743290338
$ws.Range("Q14").Value = 7105.431668887
$ws.Range("R14").Value = 63948.885019983
$ws.Range("S14").Value = 0.1805618577985511
$ws.Range("T14").Value = 0.180561857798551

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Dag1"
$ws.Range("C15").Value = "Lama4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 89.42491666666666
$ws.Range("H15").Value = 268.27475
$ws.Range("I15").Value = 0.5826613975855819
$ws.Range("J15").Value = 0.5826613975855818
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 164.182683
$ws.Range("N15").Value = 492.548049
$ws.Range("O15").Value = 0.6403319819152736
$ws.Range("P15").Value = 0.6403319819152735
$ws.Range("Q15").Value = 14682.02274538475
$ws.Range("R15").Value = 132138.2047084627
$ws.Range("S15").Value = 0.3730967275014989
$ws.Range("T15").Value = 0.3730967275014987

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Dag1"
$ws.Range("C16").Value = "Lama4"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 89.42491666666666
$ws.Range("H16").Value = 268.27475
$ws.Range("I16").Value = 0.5826613975855819
$ws.Range("J16").Value = 0.5826613975855818
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.07725
$ws.Range("N16").Value = 0.23175
$ws.Range("O16").Value = 0.0003012841835637129
$ws.Range("P16").Value = 0.0003012841835637129
$ws.Range("Q16").Value = 6.9080748125
$ws.Range("R16").Value = 62.1726733125
$ws.Range("S16").Value = 0.000175546663465664
$ws.Range("T16").Value = 0.0001755466634656639

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Dag1"
$ws.Range("C17").Value = "Lama4"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 89.42491666666666
$ws.Range("H17").Value = 268.27475
$ws.Range("I17").Value = 0.5826613975855819
$ws.Range("J17").Value = 0.5826613975855818
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 12.68555166666667
$ws.Range("N17").Value = 38.056655
$ws.Range("O17").Value = 0.04947515957212899
$ws.Range("P17").Value = 0.04947515957212898
$ws.Range("Q17").Value = 1134.404400662361
$ws.Range("R17").Value = 10209.63960596125
$ws.Range("S17").Value = 0.02882726562206636
$ws.Range("T17").Value = 0.02882726562206635
